$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the blog card in I7: ser 106 -> ser 110
$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 110"

# Selecting the edited cell updates the sheet's stored selection (C7 -> I7)
$ws.Range("I7").Select()
